# Update crypto price/volume table to reflect latest scraped values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.235.08'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.16%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.857.97'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.38%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9998'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7108'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.04%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '237.84'
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.08178'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +10.09%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3041'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.42%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.27'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.61%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08194'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.44%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.916.71'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.48%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.177'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.81%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.7083'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.77%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '89.51'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.33%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '29.268.14'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.40%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000007923'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.27%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.795'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.13%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.35'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.79%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '237.71'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.54%  '
$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.000'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.04%  '
$ws.Range("B22").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.115.84'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.40%  '
$ws.Range("E23").Value = '  -0.02%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.404'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.63%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '162.77'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.25%  '
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.953'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.73%  '
$ws.Range("B27").Value = 'Stellar'
$ws.Range("C27").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1455'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.05%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.09'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.18%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.957'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.08%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.429'
$ws.Range("D30").Style = "Normal"
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.486'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.88%  '
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.401'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.69%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.021'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.20%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05228'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.42%  '
$ws.Range("E35").Value = '  -1.66%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7079'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.03%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9992'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.44%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.676'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.56%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01856'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.73%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.730'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.90%  '
$ws.Range("B41").Value = 'Maker'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.139.89'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +6.70%  '
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9215'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.26%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4287'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.57%  '
$ws.Range("E44").Value = '  -3.32%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '70.31'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.41%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.9996'
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '102.21'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.49%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.775'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.49%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.011.02'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.08%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.184'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.67%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.978'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.08%  '
